# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 313.5
$ws.Range("J2").Value = 340.69232
$ws.Range("L2").Value = 340.69232
$ws.Range("N2").Value = -566.69232

$ws.Range("H17").Value = 3007913
$ws.Range("J17").Value = 3007913
$ws.Range("L17").Value = 9023739
$ws.Range("N17").Value = -9024075

$ws.Range("H19").Value = 1690.7333
$ws.Range("I19").Value = 849.8
$ws.Range("J19").Value = 2111.2
$ws.Range("K19").Value = 849.8
$ws.Range("L19").Value = 2111.2
$ws.Range("M19").Value = -674.8
$ws.Range("N19").Value = -2461.2

$ws.Range("H76").Value = 6409.5454
$ws.Range("I76").Value = 6325.125
$ws.Range("K76").Value = 6325.125
$ws.Range("M76").Value = -6010.125

$ws.Range("H79").Value = 6409.5454
$ws.Range("I79").Value = 6325.125
$ws.Range("K79").Value = 6325.125
$ws.Range("M79").Value = -5233.125

$ws.Range("H101").Value = 1451
$ws.Range("J101").Value = 759
$ws.Range("L101").Value = 2277
$ws.Range("N101").Value = -5521

$ws.Range("H112").Value = 5051.316
$ws.Range("J112").Value = 5370
$ws.Range("L112").Value = 16110
$ws.Range("N112").Value = -18326

$ws.Range("H137").Value = 10711.695
$ws.Range("I137").Value = 4815.269
$ws.Range("J137").Value = 15357.363
$ws.Range("K137").Value = 14445.807
$ws.Range("L137").Value = 46072.089
$ws.Range("M137").Value = -11895.807
$ws.Range("N137").Value = -51172.089

$ws.Range("H138").Value = 9630.434999999999
$ws.Range("I138").Value = 7766
$ws.Range("J138").Value = 9989.843999999999
$ws.Range("K138").Value = 23298
$ws.Range("L138").Value = 29969.532
$ws.Range("M138").Value = -18158
$ws.Range("N138").Value = -40249.532

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 964
$ws.Range("I4").Value = 931.3333
$ws.Range("J4").Value = 996.6667
$ws.Range("K4").Value = 931.3333
$ws.Range("L4").Value = 996.6667
$ws.Range("M4").Value = -815.3333
$ws.Range("N4").Value = -1228.6667

$ws.Range("H5").Value = 1725.5
$ws.Range("I5").Value = 1725.5
$ws.Range("K5").Value = 1725.5
$ws.Range("M5").Value = -1613.5

$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H11").Value = 8915
$ws.Range("J11").Value = 8915
$ws.Range("L11").Value = 8915
$ws.Range("N11").Value = -9203

$ws.Range("H12").Value = 4666.6665
$ws.Range("I12").Value = 3000
$ws.Range("J12").Value = 5500
$ws.Range("K12").Value = 3000
$ws.Range("L12").Value = 5500
$ws.Range("M12").Value = -2827
$ws.Range("N12").Value = -5846

$ws.Range("H13").Value = 500000
$ws.Range("I13").Value = 500000
$ws.Range("K13").Value = 500000
$ws.Range("M13").Value = -499856

$ws.Range("H32").Value = 71335.84
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 71335.84
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 71335.84
$ws.Range("N32").Value = -71909.84
$ws.Range("M32").ClearContents()

$ws.Range("H74").Value = 15935.939
$ws.Range("I74").Value = 3825.5557
$ws.Range("J74").Value = 30468.4
$ws.Range("K74").Value = 3825.5557
$ws.Range("L74").Value = 30468.4
$ws.Range("M74").Value = -2951.5557
$ws.Range("N74").Value = -32216.4

$ws.Range("H77").Value = 15935.939
$ws.Range("I77").Value = 3825.5557
$ws.Range("J77").Value = 30468.4
$ws.Range("K77").Value = 19127.7785
$ws.Range("L77").Value = 152342
$ws.Range("M77").Value = -14759.7785
$ws.Range("N77").Value = -161078

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1725.5
$ws.Range("I4").Value = 1725.5
$ws.Range("K4").Value = 1725.5
$ws.Range("M4").Value = -1610.5

$ws.Range("H26").Value = 35000
$ws.Range("I26").Value = 35000
$ws.Range("K26").Value = 35000
$ws.Range("M26").Value = -34708

$ws.Range("H81").Value = 69682.875
$ws.Range("J81").Value = 69682.875
$ws.Range("L81").Value = 69682.875
$ws.Range("N81").Value = -71804.875

$ws.Range("H84").Value = 69682.875
$ws.Range("J84").Value = 69682.875
$ws.Range("L84").Value = 209048.625
$ws.Range("N84").Value = -219656.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 711.84375
$ws.Range("I7").Value = 574.5
$ws.Range("J7").Value = 1307
$ws.Range("K7").Value = 574.5
$ws.Range("L7").Value = 1307
$ws.Range("M7").Value = -461.5
$ws.Range("N7").Value = -1533

$ws.Range("H59").Value = 42798.5
$ws.Range("J59").Value = 42798.5
$ws.Range("L59").Value = 42798.5
$ws.Range("N59").Value = -45088.5

$ws.Range("H62").Value = 5247.6875
$ws.Range("J62").Value = 5213.7144
$ws.Range("L62").Value = 5213.7144
$ws.Range("N62").Value = -6461.7144

$ws.Range("H65").Value = 5247.6875
$ws.Range("J65").Value = 5213.7144
$ws.Range("L65").Value = 26068.572
$ws.Range("N65").Value = -32308.572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2990
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H131").Value = 1458.54
$ws.Range("I131").Value = 971.75
$ws.Range("J131").Value = 1478.8229
$ws.Range("K131").Value = 2915.25
$ws.Range("L131").Value = 4436.468699999999
$ws.Range("M131").Value = 2124.75
$ws.Range("N131").Value = -14516.4687

$ws.Range("H137").Value = 4440.533
$ws.Range("I137").Value = 3860.8
$ws.Range("K137").Value = 11582.4
$ws.Range("M137").Value = -6482.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13903.619
$ws.Range("I80").Value = 7616.0835
$ws.Range("J80").Value = 22287
$ws.Range("K80").Value = 7616.0835
$ws.Range("L80").Value = 22287
$ws.Range("M80").Value = -6618.0835
$ws.Range("N80").Value = -24283

$ws.Range("H83").Value = 13903.619
$ws.Range("I83").Value = 7616.0835
$ws.Range("J83").Value = 22287
$ws.Range("K83").Value = 38080.4175
$ws.Range("L83").Value = 111435
$ws.Range("M83").Value = -33088.4175
$ws.Range("N83").Value = -121419

$ws.Range("H133").Value = 97628.28999999999
$ws.Range("J133").Value = 97628.28999999999
$ws.Range("L133").Value = 97628.28999999999
$ws.Range("N133").Value = -107748.29

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8105.6206
$ws.Range("I40").Value = 4056
$ws.Range("J40").Value = 15799.9
$ws.Range("K40").Value = 4056
$ws.Range("L40").Value = 15799.9
$ws.Range("M40").Value = -3920
$ws.Range("N40").Value = -16071.9

$ws.Range("H61").Value = 8250
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 8250
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 8250
$ws.Range("N61").Value = -8654
$ws.Range("M61").ClearContents()

$ws.Range("H113").Value = 8250
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 8250
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 8250
$ws.Range("N113").Value = -12590
$ws.Range("M113").ClearContents()

$ws.Range("H122").Value = 9233.083000000001
$ws.Range("J122").Value = 13333.111
$ws.Range("L122").Value = 39999.333
$ws.Range("N122").Value = -44899.333

$ws.Range("H136").Value = 14981.318
$ws.Range("I136").Value = 12673.296
$ws.Range("K136").Value = 38019.888
$ws.Range("M136").Value = -35469.888

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 39000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H136").Value = 8496.543
$ws.Range("I136").Value = 2128.739
$ws.Range("K136").Value = 6386.217000000001
$ws.Range("M136").Value = -3836.217000000001
